$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").EntireColumn.Delete()
$ws.Range("B12").Select()
